# Genetics.xlsx — "update genetics denied, logout provider"
#
# The test-data row on Sheet1 (row 2) carries a case id in column A.
# This commit rolls that id forward to a fresh value, simulating the
# provider session being logged out / the case being re-denied with a
# freshly generated case id.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2, column A ("id") -> new case id
$ws.Cells.Item(2, 1).Value = "CA-I5OS3OUL"
